$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'29.835.43"
$ws.Range('E2').Value = "'  -0.32%  "
$ws.Range('D3').Value = "'1.887.22"
$ws.Range('E3').Value = "'  -0.42%  "
$ws.Range('D4').Value = "'1.001"
$ws.Range('E4').Value = "'  +0.05%  "
$ws.Range('D5').Value = "'0.7498"
$ws.Range('E5').Value = "'  -3.33%  "
$ws.Range('D6').Value = "'242.11"
$ws.Range('E6').Value = "'  -0.74%  "
$ws.Range('D7').Value = "'1.001"
$ws.Range('E7').Value = "'  +0.03%  "
$ws.Range('D8').Value = "'0.3125"
$ws.Range('E8').Value = "'  -0.18%  "
$ws.Range('D9').Value = "'25.25"
$ws.Range('E9').Value = "'  -2.14%  "
$ws.Range('D10').Value = "'0.07108"
$ws.Range('E10').Value = "'  -3.57%  "
$ws.Range('D11').Value = "'0.08509"
$ws.Range('E11').Value = "'  +5.39%  "
$ws.Range('D12').Value = "'0.7584"
$ws.Range('E12').Value = "'  -1.96%  "
$ws.Range('D13').Value = "'1.889.50"
$ws.Range('E13').Value = "'  -0.38%  "
$ws.Range('D14').Value = "'5.364"
$ws.Range('E14').Value = "'  -2.65%  "
$ws.Range('D15').Value = "'93.24"
$ws.Range('E15').Value = "'  -1.18%  "
$ws.Range('D16').Value = "'6.126"
$ws.Range('E16').Value = "'  -1.82%  "
$ws.Range('D17').Value = "'29.797.81"
$ws.Range('E17').Value = "'  -0.73%  "
$ws.Range('D18').Value = "'13.70"
$ws.Range('E18').Value = "'  -2.21%  "
$ws.Range('D19').Value = "'242.86"
$ws.Range('E19').Value = "'  -1.90%  "
$ws.Range('D20').Value = "'0.000007827"
$ws.Range('E20').Value = "'  -0.26%  "
$ws.Range('B21').Value = "'WrappedliquidstakedEther2.0"
$ws.Range('C21').Value = "'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range('D21').Value = "'2.144.42"
$ws.Range('E21').Value = "'  -2.45%  "
$ws.Range('B22').Value = "'Dai"
$ws.Range('C22').Value = "'https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range('D22').Value = "'0.9999"
$ws.Range('E22').Value = "'  -0.06%  "
$ws.Range('D23').Value = "'7.968"
$ws.Range('E23').Value = "'  -2.58%  "
$ws.Range('E24').Value = "'  +0.06%  "
$ws.Range('D25').Value = "'0.1586"
$ws.Range('E25').Value = "'  +0.29%  "
$ws.Range('D26').Value = "'9.358"
$ws.Range('E26').Value = "'  -1.23%  "
$ws.Range('D27').Value = "'163.08"
$ws.Range('E27').Value = "'  -0.14%  "
$ws.Range('D28').Value = "'18.70"
$ws.Range('E28').Value = "'  -0.29%  "
$ws.Range('D29').Value = "'2.028"
$ws.Range('E29').Value = "'  -0.23%  "
$ws.Range('D30').Value = "'1.467"
$ws.Range('E30').Value = "'  +2.48%  "
$ws.Range('D31').Value = "'1.531"
$ws.Range('E31').Value = "'  -0.78%  "
$ws.Range('D32').Value = "'4.503"
$ws.Range('E32').Value = "'  +0.61%  "
$ws.Range('D33').Value = "'4.156"
$ws.Range('E33').Value = "'  +2.14%  "
$ws.Range('D34').Value = "'0.05415"
$ws.Range('E34').Value = "'  -2.90%  "
$ws.Range('D35').Value = "'1.238"
$ws.Range('E35').Value = "'  -0.42%  "
$ws.Range('D36').Value = "'0.7514"
$ws.Range('E36').Value = "'  -0.61%  "
$ws.Range('D37').Value = "'1.003"
$ws.Range('E37').Value = "'  -0.31%  "
$ws.Range('D38').Value = "'2.711"
$ws.Range('E38').Value = "'  +1.00%  "
$ws.Range('D39').Value = "'0.01942"
$ws.Range('E39').Value = "'  +0.46%  "
$ws.Range('E40').Value = "'  -0.66%  "
$ws.Range('D41').Value = "'0.4462"
$ws.Range('E41').Value = "'  -0.36%  "
$ws.Range('B42').Value = "'FraxShare"
$ws.Range('C42').Value = "'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range('D42').Value = "'6.104"
$ws.Range('E42').Value = "'  +1.62%  "
$ws.Range('B43').Value = "'Maker"
$ws.Range('C43').Value = "'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range('D43').Value = "'1.100.39"
$ws.Range('E43').Value = "'  -0.80%  "
$ws.Range('D44').Value = "'72.40"
$ws.Range('E44').Value = "'  -2.79%  "
$ws.Range('D45').Value = "'0.8591"
$ws.Range('E45').Value = "'  +0.84%  "
$ws.Range('E46').Value = "'  +0.05%  "
$ws.Range('D47').Value = "'7.708"
$ws.Range('E47').Value = "'  +2.26%  "
$ws.Range('D48').Value = "'102.32"
$ws.Range('E48').Value = "'  -0.34%  "
$ws.Range('D49').Value = "'1.857"
$ws.Range('E49').Value = "'  -2.12%  "
$ws.Range('D50').Value = "'3.035"
$ws.Range('E50').Value = "'  +0.86%  "
$ws.Range('D51').Value = "'2.037.31"
$ws.Range('E51').Value = "'  -0.97%  "
